$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Remove the old "A 24922-2019" row (row 4) - its data will be
#    re-inserted (with updated figures) as the new first data row.
# ------------------------------------------------------------------
$ws.Rows.Item(4).Delete()

# ------------------------------------------------------------------
# 2. Insert a brand-new row above the current row 2, pushing
#    "A 32165-2019" (old row 2) and "A 26207-2019" (old row 3) down
#    to rows 3 and 4 respectively.
# ------------------------------------------------------------------
$ws.Rows.Item(2).Insert()

# ------------------------------------------------------------------
# 3. Populate the freshly inserted row 2 with the refreshed
#    "A 24922-2019" record.
# ------------------------------------------------------------------
$ws.Range("A2").Value = "A 24922-2019"

$ws.Range("B2").Value = 43599
$ws.Range("B2").NumberFormat = "YYYY-MM-DD"

$ws.Range("C2").Value = 45179
$ws.Range("C2").NumberFormat = "YYYY-MM-DD"

$ws.Range("D2").Value = "VÄSTERBOTTENS LÄN"
$ws.Range("E2").Value = "VÄNNÄS"

$ws.Range("G2").Value = 4
$ws.Range("H2").Value = 2
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 6
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 9
$ws.Range("P2").Value = 3
$ws.Range("Q2").Value = 9

$ws.Range("R2").Value = "Blackticka`r`nRynkskinn`r`nUlltickeporing`r`nGammelgransskål`r`nGarnlav`r`nGranticka`r`nJärpe`r`nTretåig hackspett`r`nUllticka"
$ws.Range("R2").WrapText = $true

$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_VANNAS/artfynd/A 24922-2019.xlsx")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_VANNAS/kartor/A 24922-2019.png")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_VANNAS/klagomål/A 24922-2019.docx")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_VANNAS/klagomålsmail/A 24922-2019.docx")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_VANNAS/tillsyn/A 24922-2019.docx")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_VANNAS/tillsynsmail/A 24922-2019.docx")'

# Match the fixed row height used throughout the sheet.
$ws.Rows.Item(2).RowHeight = 15

# ------------------------------------------------------------------
# 4. Every data row's "Förändrad" (column C) date moved one day
#    forward, from 45178 to 45179. Apply this uniformly to all data
#    rows (2 through 171) - this also covers the just-inserted row 2.
# ------------------------------------------------------------------
for ($r = 2; $r -le 171; $r++) {
    $cell = $ws.Range("C$r")
    if ($cell.Value() -ne $null) {
        $cell.Value = 45179
        $cell.NumberFormat = "YYYY-MM-DD"
    }
}
